# Update recomputed NATMI ligand-receptor metrics (new TPM input) for rows 2-10
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1346003333333333
$ws.Range("H2").Value = 0.403801
$ws.Range("I2").Value = 0.009651054304565105
$ws.Range("J2").Value = 0.009651054304565105
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.01848533333333334
$ws.Range("N2").Value = 0.05545600000000001
$ws.Range("O2").Value = 0.001625201930372746
$ws.Range("P2").Value = 0.001625201930372746
$ws.Range("Q2").Value = 0.002488132028444445
$ws.Range("R2").Value = 0.022393188256
$ws.Range("S2").Value = 0.00001568491208591141
$ws.Range("T2").Value = 0.0000156849120859114

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1346003333333333
$ws.Range("H3").Value = 0.403801
$ws.Range("I3").Value = 0.009651054304565105
$ws.Range("J3").Value = 0.009651054304565105
$ws.Range("O3").Value = 0.002698334581238102
$ws.Range("P3").Value = 0.002698334581238102
$ws.Range("Q3").Value = 0.004131063697111112
$ws.Range("R3").Value = 0.037179573274
$ws.Range("S3").Value = 0.00002604177357541487
$ws.Range("T3").Value = 0.00002604177357541486

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1346003333333333
$ws.Range("H4").Value = 0.403801
$ws.Range("I4").Value = 0.009651054304565105
$ws.Range("J4").Value = 0.009651054304565105
$ws.Range("M4").Value = 11.32499966666667
$ws.Range("N4").Value = 33.974999
$ws.Range("O4").Value = 0.9956764634883892
$ws.Range("P4").Value = 0.995676463488389
$ws.Range("Q4").Value = 1.524348730133223
$ws.Range("R4").Value = 13.719138571199
$ws.Range("S4").Value = 0.00960932761890378
$ws.Range("T4").Value = 0.009609327618903778

# Row 5
$ws.Range("I5").Value = 0.8124788779145131
$ws.Range("J5").Value = 0.8124788779145132
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.01848533333333334
$ws.Range("N5").Value = 0.05545600000000001
$ws.Range("O5").Value = 0.001625201930372746
$ws.Range("P5").Value = 0.001625201930372746
$ws.Range("Q5").Value = 0.2094646506773334
$ws.Range("R5").Value = 1.885181856096
$ws.Range("S5").Value = 0.001320442240773749
$ws.Range("T5").Value = 0.001320442240773749

# Row 6
$ws.Range("I6").Value = 0.8124788779145131
$ws.Range("J6").Value = 0.8124788779145132
$ws.Range("O6").Value = 0.002698334581238102
$ws.Range("P6").Value = 0.002698334581238102
$ws.Range("S6").Value = 0.002192339852802261
$ws.Range("T6").Value = 0.002192339852802261

# Row 7
$ws.Range("I7").Value = 0.8124788779145131
$ws.Range("J7").Value = 0.8124788779145132
$ws.Range("M7").Value = 11.32499966666667
$ws.Range("N7").Value = 33.974999
$ws.Range("O7").Value = 0.9956764634883892
$ws.Range("P7").Value = 0.995676463488389
$ws.Range("Q7").Value = 128.3280672478677
$ws.Range("R7").Value = 1154.952605230809
$ws.Range("S7").Value = 0.8089660958209371
$ws.Range("T7").Value = 0.8089660958209371

# Row 8
$ws.Range("G8").Value = 2.4807
$ws.Range("H8").Value = 7.4421
$ws.Range("I8").Value = 0.1778700677809217
$ws.Range("J8").Value = 0.1778700677809217
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.01848533333333334
$ws.Range("N8").Value = 0.05545600000000001
$ws.Range("O8").Value = 0.001625201930372746
$ws.Range("P8").Value = 0.001625201930372746
$ws.Range("Q8").Value = 0.04585656640000001
$ws.Range("R8").Value = 0.4127090976
$ws.Range("S8").Value = 0.0002890747775130851
$ws.Range("T8").Value = 0.0002890747775130851

# Row 9
$ws.Range("G9").Value = 2.4807
$ws.Range("H9").Value = 7.4421
$ws.Range("I9").Value = 0.1778700677809217
$ws.Range("J9").Value = 0.1778700677809217
$ws.Range("O9").Value = 0.002698334581238102
$ws.Range("P9").Value = 0.002698334581238102
$ws.Range("Q9").Value = 0.07613599060000001
$ws.Range("R9").Value = 0.6852239154
$ws.Range("S9").Value = 0.0004799529548604262
$ws.Range("T9").Value = 0.0004799529548604261

# Row 10
$ws.Range("G10").Value = 2.4807
$ws.Range("H10").Value = 7.4421
$ws.Range("I10").Value = 0.1778700677809217
$ws.Range("J10").Value = 0.1778700677809217
$ws.Range("M10").Value = 11.32499966666667
$ws.Range("N10").Value = 33.974999
$ws.Range("O10").Value = 0.9956764634883892
$ws.Range("P10").Value = 0.995676463488389
$ws.Range("Q10").Value = 28.09392667310001
$ws.Range("R10").Value = 252.8453400579
$ws.Range("S10").Value = 0.1771010400485482
$ws.Range("T10").Value = 0.1771010400485482

